# data : case 1
# Update the two data columns on the active sheet with the new computed
# values, and nudge the column widths to match the new (slightly wider)
# layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column widths (closest values reachable through the ColumnWidth property;
# target stored widths are 15.7109375 and 16.42578125 character units).
$ws.Columns.Item(1).ColumnWidth = 14.833333333333334
$ws.Columns.Item(2).ColumnWidth = 15.666666666666666

$colA = @(
  -0.1026037900120258,
  -0.080198531312594135,
  -0.029283859676075608,
  -0.02097351071272513,
  -0.017700902759337112,
  0.00080702943691690621,
  0.010987517476769071,
  0.021022308505631138,
  0.02307785086382852,
  0.025123961409518003,
  0.028124287461001707,
  -0.0094290460888473149,
  -0.0058175800422350932,
  -0.0048216673550269817,
  -0.0038206864027463538,
  -0.0018204646698585236,
  0.00018533356529903955,
  -0.016100207346429585,
  -0.012090675023090913,
  -0.0080157483590532763,
  -0.0040056197244613401,
  -0.045712831275540822,
  -0.040499887116864031,
  -0.020098991127214028,
  -0.015592174490896227,
  -0.028210665343966213,
  -0.025621198395235911,
  -0.023100165281187302,
  -0.01576215059408792,
  0.044324760246357453,
  -0.014021159071152312,
  -0.0040008422288817513
)

$colB = @(
  0.10229873858367,
  0.078996690892987331,
  0.028973510603316655,
  0.020700902710778735,
  0.016785145016782188,
  -0.00098751761136384175,
  -0.011022308638140466,
  -0.021077850891420002,
  -0.023123961430615569,
  -0.025124287493454744,
  -0.028129889994789714,
  0.0093175800037075263,
  0.0057295487839956749,
  0.0048206863950879253,
  0.0038204646497970174,
  0.001814666414702959,
  -0.00019055030630177328,
  0.016090674972449648,
  0.012015748304308183,
  0.008005619669102515,
  0.0039999999442095202,
  0.045499887041788867,
  0.040098990859671169,
  0.019999999728598894,
  0.015564528403443134,
  0.028121198350101295,
  0.025100165236158656,
  0.022762150483909949,
  0.015675238987268969,
  -0.044611959510094668,
  0.014000842091242305,
  0.003999999936487697
)

for ($i = 0; $i -lt $colA.Length; $i++) {
  $row = $i + 1
  $ws.Cells.Item($row, 1).Value = $colA[$i]
  $ws.Cells.Item($row, 2).Value = $colB[$i]
}
